# Update "想去人数" (want-to-go count, column F) values on the
# "展览" and "全部类型" sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> list of (row, newValue)
$updates = @{
    "展览"   = @(
        @{ Row = 2;  Value = 569 },
        @{ Row = 4;  Value = 1239 },
        @{ Row = 6;  Value = 14021 },
        @{ Row = 7;  Value = 15341 },
        @{ Row = 20; Value = 1185 },
        @{ Row = 23; Value = 5963 },
        @{ Row = 24; Value = 955 },
        @{ Row = 26; Value = 5512 },
        @{ Row = 29; Value = 96 },
        @{ Row = 30; Value = 436 }
    )
    "全部类型" = @(
        @{ Row = 3;  Value = 569 },
        @{ Row = 5;  Value = 1239 },
        @{ Row = 7;  Value = 14021 },
        @{ Row = 8;  Value = 15341 },
        @{ Row = 21; Value = 1185 },
        @{ Row = 25; Value = 5963 },
        @{ Row = 26; Value = 955 },
        @{ Row = 28; Value = 5512 },
        @{ Row = 31; Value = 96 },
        @{ Row = 32; Value = 436 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $cell = $ws.Cells.Item($entry.Row, 6)  # Column F = 6
        $cell.Value = $entry.Value
    }
}
